$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data row with the newer COVID-19 snapshot (20200321 pull) ---
$ws.Range("A2").Value = 43911   # 2020-03-21
$ws.Range("B2").Value = 5018
$ws.Range("C2").Value = 1035
$ws.Range("D2").Value = 233
$ws.Range("E2").Value = 4257
$ws.Range("F2").Value = 373
$ws.Range("G2").Value = 280
$ws.Range("H2").Value = 108

# --- Strip the old thin-border / fill / custom numFmt formatting from the row ---
# Counts (B2:H2) revert to the plain default style.
$ws.Range("B2:H2").Style = "Normal"

# Date cell keeps a date format but also loses its border/center alignment,
# reusing the built-in short-date number format (numFmtId 14).
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").NumberFormat = "mm-dd-yy"

# --- Move the active selection ---
$ws.Range("C4").Select() | Out-Null
